# Apply the "Revision hasta TIPO DE PERSONA TITULAR" corrections to the
# Sotara / Medellin / San Antero / San Juan del Cesar / Santander rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-87: "SOTARA PAISPAMBA" -> "SOTARÁ" (col D) / "SOTARA" (col F)
for ($r = 2; $r -le 87; $r++) {
    $ws.Cells.Item($r, 4).Value = "SOTARÁ"
    $ws.Cells.Item($r, 6).Value = "SOTARA"
}

# Row 88: "BOLIVAR" -> "SAN ANTERO" (cols D and F); C/E (CORDOBA) unchanged
$ws.Cells.Item(88, 4).Value = "SAN ANTERO"
$ws.Cells.Item(88, 6).Value = "SAN ANTERO"

# Row 89: was entirely empty in C:F -> now populated
$ws.Cells.Item(89, 3).Value = "ANTIOQUIA"
$ws.Cells.Item(89, 4).Value = "MEDELLÍN"
$ws.Cells.Item(89, 5).Value = "ANTIOQUIA"
$ws.Cells.Item(89, 6).Value = "MEDELLIN"

# Row 90: "ANTIOQUIA" -> "SANTANDER" (cols C and E); D/F (PUERTO WILCHES) unchanged
$ws.Cells.Item(90, 3).Value = "SANTANDER"
$ws.Cells.Item(90, 5).Value = "SANTANDER"

# Row 91: "SAN JUAN" -> "SAN JUAN DEL CESAR" (cols D and F); C/E (LA GUAJIRA) unchanged
$ws.Cells.Item(91, 4).Value = "SAN JUAN DEL CESAR"
$ws.Cells.Item(91, 6).Value = "SAN JUAN DEL CESAR"

# Rows 92-103: were entirely empty in C:F -> now populated
for ($r = 92; $r -le 103; $r++) {
    $ws.Cells.Item($r, 3).Value = "ANTIOQUIA"
    $ws.Cells.Item($r, 4).Value = "MEDELLÍN"
    $ws.Cells.Item($r, 5).Value = "ANTIOQUIA"
    $ws.Cells.Item($r, 6).Value = "MEDELLIN"
}
